$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.098.64'
$ws.Cells.Item(2, 5).Value = '  +0.16%  '
$ws.Cells.Item(3, 4).Value = '1.780.34'
$ws.Cells.Item(3, 5).Value = '  -0.40%  '
$ws.Cells.Item(4, 5).Value = '  +0.25%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '225.51'
$ws.Cells.Item(5, 5).Value = '  -0.47%  '
$ws.Cells.Item(6, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 5).Value = '  +0.24%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '31.81'
$ws.Cells.Item(8, 5).Value = '  -1.02%  '
$ws.Cells.Item(9, 5).Value = '  -1.35%  '
$ws.Cells.Item(10, 5).Value = '  +0.32%  '
$ws.Cells.Item(11, 5).Value = '  +0.78%  '
$ws.Cells.Item(12, 4).Value = '2.036.82'
$ws.Cells.Item(12, 5).Value = '  -0.44%  '
$ws.Cells.Item(13, 4).Value = '1.785.91'
$ws.Cells.Item(13, 5).Value = '  -0.35%  '
$ws.Cells.Item(14, 5).Value = '  -3.20%  '
$ws.Cells.Item(15, 4).Value = '34.089.47'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.622'
$ws.Cells.Item(16, 5).Value = '  +0.33%  '
$ws.Cells.Item(17, 5).Value = '  -0.02%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '67.62'
$ws.Cells.Item(18, 5).Value = '  -0.15%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '244.77'
$ws.Cells.Item(19, 5).Value = '  +0.92%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0786'
$ws.Cells.Item(20, 5).Value = '  +1.68%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.91'
$ws.Cells.Item(21, 5).Value = '  +1.98%  '
$ws.Cells.Item(22, 5).Value = '  +0.28%  '
$ws.Cells.Item(23, 5).Value = '  +0.35%  '
$ws.Cells.Item(24, 5).Value = '  -1.03%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '161.37'
$ws.Cells.Item(25, 5).Value = '  -0.34%  '
$ws.Cells.Item(26, 5).Value = '  -0.50%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '16.24'
$ws.Cells.Item(27, 5).Value = '  +0.24%  '
$ws.Cells.Item(28, 5).Value = '  +1.02%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  +0.31%  '
$ws.Cells.Item(30, 5).Value = '  -0.11%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.0517'
$ws.Cells.Item(31, 5).Value = '  +0.30%  '
$ws.Cells.Item(32, 5).Value = '  +1.84%  '
$ws.Cells.Item(33, 5).Value = '  +2.89%  '
$ws.Cells.Item(34, 5).Value = '  -2.36%  '
$ws.Cells.Item(35, 4).Value = '1.448.53'
$ws.Cells.Item(35, 5).Value = '  +3.67%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.46'
$ws.Cells.Item(36, 5).Value = '  +4.81%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.653'
$ws.Cells.Item(37, 5).Value = '  +0.29%  '
$ws.Cells.Item(38, 5).Value = '  +1.16%  '
$ws.Cells.Item(39, 5).Value = '  -0.29%  '
$ws.Cells.Item(40, 5).Value = '  +1.50%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '80.36'
$ws.Cells.Item(41, 5).Value = '  +0.28%  '
$ws.Cells.Item(42, 5).Value = '  +1.38%  '
$ws.Cells.Item(43, 5).Value = '  -0.29%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '13.69'
$ws.Cells.Item(44, 5).Value = '  -0.11%  '
$ws.Cells.Item(45, 5).Value = '  +1.69%  '
$ws.Cells.Item(46, 5).Value = '  -0.08%  '
$ws.Cells.Item(47, 5).Value = '  -0.21%  '
$ws.Cells.Item(48, 5).Value = '  -0.39%  '
$ws.Cells.Item(49, 2).Value = 'Quant'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '104.24'
$ws.Cells.Item(49, 5).Value = '  -3.16%  '
$ws.Cells.Item(50, 2).Value = 'PaxDollar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.00'
$ws.Cells.Item(50, 5).Value = '  +0.32%  '
$ws.Cells.Item(51, 4).Value = '0.0₆0129'
$ws.Cells.Item(51, 5).Value = '  -7.01%  '
